# Aggiornamento fino a 6 gennaio 2022
# Appends new daily rows (465-491) to Sheet1, mirroring the formatting of
# the last existing data row (464) for column A (date style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  @(465, 44539, 4, 18, 562.6758361988121),
  @(466, 44540, 1, 16, 500.1562988433885),
  @(467, 44541, 0, 16, 500.1562988433885),
  @(468, 44542, 3, 11, 343.8574554548296),
  @(469, 44543, 2, 13, 406.3769928102532),
  @(470, 44544, 0, 11, 343.8574554548296),
  @(471, 44545, 0, 10, 312.5976867771179),
  @(472, 44546, 0, 6, 187.5586120662707),
  @(473, 44547, 0, 5, 156.2988433885589),
  @(474, 44548, 2, 7, 218.8183807439825),
  @(475, 44550, 0, 4, 125.0390747108471),
  @(476, 44551, 2, 4, 125.0390747108471),
  @(477, 44552, 0, 4, 125.0390747108471),
  @(478, 44553, 1, 5, 156.2988433885589),
  @(479, 44554, 0, 5, 156.2988433885589),
  @(480, 44555, 7, 12, 375.1172241325414),
  @(481, 44556, 4, 14, 437.636761487965),
  @(482, 44557, 4, 18, 562.6758361988121),
  @(483, 44558, 10, 26, 812.7539856205065),
  @(484, 44559, 1, 27, 844.0137542982183),
  @(485, 44560, 3, 29, 906.5332916536416),
  @(486, 44561, 8, 37, 1156.611441075336),
  @(487, 44562, 1, 31, 969.0528290090654),
  @(488, 44563, 11, 38, 1187.871209753048),
  @(489, 44564, 17, 51, 1594.248202563301),
  @(490, 44565, 6, 47, 1469.209127852454),
  @(491, 44566, 9, 55, 1719.287277274148)
)

$templateCell = $ws.Range("A464")

foreach ($row in $newRows) {
    $r = $row[0]

    $templateCell.Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
